# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet
#    and fill it with the fund-holding detail rows for 2022-Q1 (same layout
#    as the other quarterly sheets: 基金代码/基金名称/基金规模/股票总仓位/
#    仓位占比/持有市值(亿元)/仓位排名).
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet (date, holding
#    count, holding value), pushing the existing rows down and renumbering
#    the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" worksheet, positioned right before "总计"
# ---------------------------------------------------------------------
# NOTE: worksheet handles in this host resolve by tab *position*, not a
# stable identity. Inserting a sheet shifts everything after it, so a
# handle obtained beforehand (e.g. "总计") silently starts pointing at the
# newly inserted sheet afterwards. Re-fetch "总计" by name once the insert
# (and any rename) has happened, instead of reusing the pre-insert handle.
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$q1.Name = "2022-Q1"

# Match the page margins used by the other data sheets (values are in points;
# 72pt = 1in).
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1

# Header row
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $q1.Cells.Item(1, $c + 2)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: index, code, name, size, total stock position, position ratio,
# holding value (100M yuan), position rank
$rows = @(
    @("501208", "中欧创新未来18个月封闭运作混合A", "87.66", "63.39", "2.87", "2.5158", 6),
    @("213001", "宝盈鸿利收益灵活配置混合A", "17.98", "90.37", "8.62", "1.5499", 2),
    @("398001", "中海优质成长混合", "14.42", "90.86", "5.87", "0.8465", 4),
    @("010751", "宝盈优质成长混合A", "5.64", "92.80", "6.68", "0.3768", 3),
    @("001749", "招商中国机遇股票", "4.57", "94.84", "5.48", "0.2504", 3),
    @("001543", "宝盈新锐灵活配置混合A", "3.21", "93.26", "7.30", "0.2343", 3),
    @("002103", "招商康泰灵活配置混合", "1.95", "39.68", "4.42", "0.0862", 2),
    @("007581", "宝盈鸿利收益灵活配置混合C", "0.73", "90.37", "8.62", "0.0629", 2),
    @("010752", "宝盈优质成长混合C", "0.78", "92.80", "6.68", "0.0521", 3),
    @("007578", "宝盈新锐灵活配置混合C", "0.20", "93.26", "7.30", "0.0146", 3),
    @("001252", "中海进取收益灵活配置混合", "0.23", "92.65", "5.02", "0.0115", 4),
    @("001849", "前海开源强势共识100强等权重股票", "0.12", "92.23", "1.06", "0.0013", 8)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $idxCell = $q1.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    for ($c = 0; $c -lt 5; $c++) {
        $q1.Cells.Item($r, $c + 2).Value = $row[$c]
    }
    $q1.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 2: prepend the 2022-Q1 summary row to the "总计" sheet
# ---------------------------------------------------------------------
$totals = @(
    @("2022-Q1", 12, 6),
    @("2021-Q4", 13, 3.26),
    @("2021-Q3", 2, 0.03),
    @("2021-Q2", 1, 0.04)
)

for ($i = 0; $i -lt $totals.Length; $i++) {
    $r = $i + 2
    $t = $totals[$i]

    $idxCell = $total.Cells.Item($r, 1)
    $idxCell.Value = $i
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $total.Cells.Item($r, 2).Value = $t[0]
    $total.Cells.Item($r, 3).Value = $t[1]
    $total.Cells.Item($r, 4).Value = $t[2]
}
